# simple-nlp.xlsx: rename Sheet1 -> ipopt, add a new "tsp" sheet with a
# TSP distance matrix, repoint the chart/defined-names at the renamed
# sheet, and drop the stale _xlchart.* helper names.

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet; this cascades into every defined name
# and solver_* reference that already pointed at "Sheet1" -----------------
$ipopt = $wb.Worksheets.Item(1)
$ipopt.Name = "ipopt"

# the chart's series formula is not covered by the sheet rename, so
# repoint it explicitly at the renamed sheet ------------------------------
$chart = $ipopt.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(,ipopt!`$B`$4:`$B`$204,ipopt!`$C`$4:`$C`$204,1)"

# drop the orphaned chart helper names, no longer needed -------------------
$wb.Names.Item("_xlchart.v1.0").Delete()
$wb.Names.Item("_xlchart.v1.1").Delete()
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# --- new "tsp" sheet, right after "ipopt" ---------------------------------
$tsp = $wb.Worksheets.Add($null, $ipopt)
$tsp.Name = "tsp"

$labels = @("C1", "C2", "C3", "C4", "C5", "C6")
$dist = @(
    @($null, 16, 63, 21, 20, 66),
    @(57, $null, 40, 46, 69, 42),
    @(23, 11, $null, 55, 53, 47),
    @(71, 53, 58, $null, 47, 5),
    @(27, 79, 53, 35, $null, 30),
    @(57, 47, 51, 17, 24, $null)
)

for ($i = 0; $i -lt 6; $i++) {
    $tsp.Cells.Item(1, $i + 2).Value = $labels[$i]   # B1:G1 header row
    $tsp.Cells.Item($i + 2, 1).Value = $labels[$i]   # A2:A7 row labels
}

for ($r = 0; $r -lt 6; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $tsp.Cells.Item($r + 2, $c + 2)
        if ($r -eq $c) {
            $cell.Formula = "=`$M`$3"
        } else {
            $cell.Value = $dist[$r][$c]
        }
    }
}

$tsp.Range("B6").Select()
